$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 78, pushing the existing rows
# 78..167 down to 80..169 (matches dimension growing from R167 to R169).
$ws.Range("A78:R79").EntireRow.Insert()

# Common (fixed) column values shared by every data row in this sheet.
$mercadoId = 8
$mercado = "Terminal La Palmera de La Serena"
$region = "Coquimbo"
$codreg = 4
$categoriaId = 100112021
$categoria = "Ají"
$clasificacion = "Hortaliza"

# New row 78: Ají / Americana (o) / Primera, Provincia de Limarí
$ws.Cells.Item(78, 1).Value = $mercadoId
$ws.Cells.Item(78, 2).Value = $mercado
$ws.Cells.Item(78, 3).Value = $region
$ws.Cells.Item(78, 4).Value = 44571
$ws.Cells.Item(78, 5).Value = $codreg
$ws.Cells.Item(78, 6).Value = $categoriaId
$ws.Cells.Item(78, 7).Value = $categoria
$ws.Cells.Item(78, 8).Value = "Americana (o)"
$ws.Cells.Item(78, 9).Value = "Primera"
$ws.Cells.Item(78, 10).Value = 600
$ws.Cells.Item(78, 11).Value = 15000
$ws.Cells.Item(78, 12).Value = 16000
$ws.Cells.Item(78, 13).Value = 15500
$ws.Cells.Item(78, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(78, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(78, 16).Value = 1033
$ws.Cells.Item(78, 17).Value = 15
$ws.Cells.Item(78, 18).Value = $clasificacion

# New row 79: Ají / Americana (o) / Segunda, Provincia de Limarí
$ws.Cells.Item(79, 1).Value = $mercadoId
$ws.Cells.Item(79, 2).Value = $mercado
$ws.Cells.Item(79, 3).Value = $region
$ws.Cells.Item(79, 4).Value = 44571
$ws.Cells.Item(79, 5).Value = $codreg
$ws.Cells.Item(79, 6).Value = $categoriaId
$ws.Cells.Item(79, 7).Value = $categoria
$ws.Cells.Item(79, 8).Value = "Americana (o)"
$ws.Cells.Item(79, 9).Value = "Segunda"
$ws.Cells.Item(79, 10).Value = 400
$ws.Cells.Item(79, 11).Value = 12000
$ws.Cells.Item(79, 12).Value = 13000
$ws.Cells.Item(79, 13).Value = 12500
$ws.Cells.Item(79, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(79, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(79, 16).Value = 833
$ws.Cells.Item(79, 17).Value = 15
$ws.Cells.Item(79, 18).Value = $clasificacion

# Match the date-format style used by the rest of column D.
$ws.Range("D78:D79").NumberFormat = $ws.Range("D80").NumberFormat
